# Fix an error in the dairy manure rate on the "gompertz" sheet.
# Row 4 (Dairy) column B held a unit-conversion constant (lb->kg, 0.453592)
# that should instead have been 0.018 (1.8E-2).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("gompertz")

# Correct the Dairy manure rate constant in B4.
$ws.Range("B4").Value = 0.018

# Recalculate so the dependent shared formulas (L4:Y4, etc.) update.
$excel.Calculate()

# Reflect the selection left behind in the source edit (row 4 highlighted).
$originalActive = $wb.ActiveSheet
$ws.Select()
$ws.Range("G4:Y4").Select()
$originalActive.Select()
